# "update with new MeSH table"
# Extends both concept tables with 2020-2024 (columns V:Z) and refreshes the
# existing 2000-2019 (columns B:U) figures with the new MeSH-derived counts.
$wb = $excel.ActiveWorkbook

$xlPasteFormats = -4122

# --- Sheet: No. of pmid ---
$ws1 = $wb.Worksheets.Item("No. of pmid")

# New year headers V1:Z1 (2020-2024). Pre-format as text so the numeric-
# looking labels are stored as strings (matching 2000-2019), then paste just
# the format from U1 so the new cells share its style instead of getting a
# separate "Text"-number-format style.
$ws1.Range("V1:Z1").NumberFormat = "@"
$newHeaders = @("2020", "2021", "2022", "2023", "2024")
for ($i = 0; $i -lt $newHeaders.Length; $i++) {
    $ws1.Cells.Item(1, 22 + $i).Value = $newHeaders[$i]
}
$ws1.Range("U1").Copy() | Out-Null
$ws1.Range("V1:Z1").PasteSpecial($xlPasteFormats) | Out-Null
$excel.CutCopyMode = $false

# Updated per-year counts for rows 2-4 (columns B:Z => years 2000-2024)
# Row 2: Artificial intelligence
$row2Vals = @(38, 36, 50, 58, 59, 54, 79, 96, 93, 102, 110, 115, 129, 130, 132, 130, 107, 130, 122, 152, 195, 249, 244, 183, 2)
for ($i = 0; $i -lt $row2Vals.Length; $i++) {
    $ws1.Cells.Item(2, 2 + $i).Value = $row2Vals[$i]
}
# Row 3: Machine learning
$row3Vals = @(4, 6, 12, 11, 13, 9, 31, 24, 24, 37, 38, 28, 33, 35, 45, 43, 35, 45, 48, 65, 80, 110, 114, 78, 0)
for ($i = 0; $i -lt $row3Vals.Length; $i++) {
    $ws1.Cells.Item(3, 2 + $i).Value = $row3Vals[$i]
}
# Row 4: Natural language processing
$row4Vals = @(0, 1, 0, 1, 0, 0, 1, 0, 1, 0, 0, 2, 1, 2, 2, 3, 1, 5, 2, 4, 1, 1, 10, 3, 0)
for ($i = 0; $i -lt $row4Vals.Length; $i++) {
    $ws1.Cells.Item(4, 2 + $i).Value = $row4Vals[$i]
}

# --- Sheet: No. of doi ---
$ws2 = $wb.Worksheets.Item("No. of doi")

# New year headers V1:Z1 (2020-2024). Pre-format as text so the numeric-
# looking labels are stored as strings (matching 2000-2019), then paste just
# the format from U1 so the new cells share its style instead of getting a
# separate "Text"-number-format style.
$ws2.Range("V1:Z1").NumberFormat = "@"
$newHeaders = @("2020", "2021", "2022", "2023", "2024")
for ($i = 0; $i -lt $newHeaders.Length; $i++) {
    $ws2.Cells.Item(1, 22 + $i).Value = $newHeaders[$i]
}
$ws2.Range("U1").Copy() | Out-Null
$ws2.Range("V1:Z1").PasteSpecial($xlPasteFormats) | Out-Null
$excel.CutCopyMode = $false

# Updated per-year counts for rows 2-4 (columns B:Z => years 2000-2024)
# Row 2: Artificial intelligence
$row2Vals = @(34, 36, 41, 58, 54, 52, 74, 91, 89, 93, 107, 108, 125, 125, 126, 124, 102, 128, 117, 151, 191, 247, 244, 183, 2)
for ($i = 0; $i -lt $row2Vals.Length; $i++) {
    $ws2.Cells.Item(2, 2 + $i).Value = $row2Vals[$i]
}
# Row 3: Machine learning
$row3Vals = @(4, 6, 12, 11, 13, 9, 29, 23, 24, 37, 36, 28, 31, 35, 45, 41, 34, 45, 47, 65, 80, 111, 114, 78, 0)
for ($i = 0; $i -lt $row3Vals.Length; $i++) {
    $ws2.Cells.Item(3, 2 + $i).Value = $row3Vals[$i]
}
# Row 4: Natural language processing
$row4Vals = @(0, 1, 0, 1, 0, 0, 1, 0, 1, 0, 0, 2, 1, 2, 2, 3, 1, 5, 2, 4, 1, 1, 10, 3, 0)
for ($i = 0; $i -lt $row4Vals.Length; $i++) {
    $ws2.Cells.Item(4, 2 + $i).Value = $row4Vals[$i]
}
